$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -188

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17302.6
$ws.Range("I21").Value = 16205.2
$ws.Range("J21").Value = 18400
$ws.Range("K21").Value = 16205.2
$ws.Range("L21").Value = 18400
$ws.Range("M21").Value = -15737.2
$ws.Range("N21").Value = -19336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 17302.6
$ws.Range("I23").Value = 16205.2
$ws.Range("J23").Value = 18400
$ws.Range("K23").Value = 16205.2
$ws.Range("L23").Value = 18400
$ws.Range("M23").Value = -15971.2
$ws.Range("N23").Value = -18868

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3298.0908
$ws.Range("I28").Value = 3199.875
$ws.Range("K28").Value = 3199.875
$ws.Range("M28").Value = -2714.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3714.1428
$ws.Range("I74").Value = 3599.8
$ws.Range("K74").Value = 3599.8
$ws.Range("M74").Value = -2663.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3714.1428
$ws.Range("I77").Value = 3599.8
$ws.Range("K77").Value = 17999
$ws.Range("M77").Value = -13319

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 125000730
$ws.Range("I127").Value = 500000260
$ws.Range("J127").Value = 883.3333
$ws.Range("K127").Value = 1500000780
$ws.Range("L127").Value = 2649.9999
$ws.Range("M127").Value = -1499995820
$ws.Range("N127").Value = -12569.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1701.8334
$ws.Range("I45").Value = 1486.0769
$ws.Range("J45").Value = 2262.8
$ws.Range("K45").Value = 1486.0769
$ws.Range("L45").Value = 2262.8
$ws.Range("M45").Value = -1109.0769
$ws.Range("N45").Value = -3016.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1292.6765
$ws.Range("I20").Value = 1025.3914
$ws.Range("J20").Value = 1851.5454
$ws.Range("K20").Value = 1025.3914
$ws.Range("L20").Value = 1851.5454
$ws.Range("M20").Value = -778.3914
$ws.Range("N20").Value = -2345.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4946.744
$ws.Range("I134").Value = 4090.5625
$ws.Range("J134").Value = 7437.4546
$ws.Range("K134").Value = 12271.6875
$ws.Range("L134").Value = 22312.3638
$ws.Range("M134").Value = -9736.6875
$ws.Range("N134").Value = -27382.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 261.3
$ws.Range("I22").Value = 142.85715
$ws.Range("K22").Value = 142.85715
$ws.Range("M22").Value = 207.14285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 9000
$ws.Range("N22").Value = -9338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 9000
$ws.Range("N27").Value = -9204

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 452.63635
$ws.Range("I34").Value = 375.75
$ws.Range("J34").Value = 496.57144
$ws.Range("K34").Value = 1127.25
$ws.Range("L34").Value = 1489.71432
$ws.Range("M34").Value = -1043.25
$ws.Range("N34").Value = -1657.71432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2030.5714
$ws.Range("J39").Value = 2030.5714
$ws.Range("L39").Value = 6091.7142
$ws.Range("N39").Value = -6679.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1466.6666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1466.6666
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 4399.9998
$ws.Range("N55").Value = -4753.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 866.03
$ws.Range("I131").Value = 548.0909
$ws.Range("J131").Value = 905.32587
$ws.Range("K131").Value = 1644.2727
$ws.Range("L131").Value = 2715.97761
$ws.Range("M131").Value = 3395.7273
$ws.Range("N131").Value = -12795.97761

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3135.1155
$ws.Range("I132").Value = 3212.5557
$ws.Range("J132").Value = 2960.875
$ws.Range("K132").Value = 9637.667099999999
$ws.Range("L132").Value = 8882.625
$ws.Range("M132").Value = -7107.667099999999
$ws.Range("N132").Value = -13942.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 502.8684
$ws.Range("I22").Value = 314.75
$ws.Range("J22").Value = 525
$ws.Range("K22").Value = 314.75
$ws.Range("L22").Value = 525
$ws.Range("M22").Value = -19.75
$ws.Range("N22").Value = -1115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 502.8684
$ws.Range("I27").Value = 314.75
$ws.Range("J27").Value = 525
$ws.Range("K27").Value = 314.75
$ws.Range("L27").Value = 525
$ws.Range("M27").Value = -207.75
$ws.Range("N27").Value = -739

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 768.5
$ws.Range("I46").Value = 829.1
$ws.Range("J46").Value = 718
$ws.Range("K46").Value = 829.1
$ws.Range("L46").Value = 718
$ws.Range("M46").Value = -641.1
$ws.Range("N46").Value = -1094

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2943.1943
$ws.Range("I132").Value = 2951
$ws.Range("J132").Value = 2939.2917
$ws.Range("K132").Value = 8853
$ws.Range("L132").Value = 8817.875100000001
$ws.Range("M132").Value = -6323
$ws.Range("N132").Value = -13877.8751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9260601
$ws.Range("I136").Value = 10753849
$ws.Range("J136").Value = 2460
$ws.Range("K136").Value = 32261547
$ws.Range("L136").Value = 7380
$ws.Range("M136").Value = -32258997
$ws.Range("N136").Value = -12480

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 351.75
$ws.Range("I113").Value = 368
$ws.Range("J113").Value = 303
$ws.Range("K113").Value = 1104
$ws.Range("L113").Value = 909
$ws.Range("M113").Value = 1066
$ws.Range("N113").Value = -5249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4169440.2
$ws.Range("I132").Value = 6063546.5
$ws.Range("J132").Value = 2406.5334
$ws.Range("K132").Value = 18190639.5
$ws.Range("L132").Value = 7219.600199999999
$ws.Range("M132").Value = -18188109.5
$ws.Range("N132").Value = -12279.6002
